$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Subject 103 (row 4): clarify that the data actually ended up under 104
$ws.Range("D4").Value = "8 presses, data under 104"

# Subject 104 (row 5): add name, note about mix-up, and scan date
$ws.Range("B5").Value = "inbal"
$ws.Range("C5").Value = "gur-arye"
$ws.Range("D5").Value = "8 presses, data un,der 103, first experimental run is bad, first auditory localizer also"

# Copy the date formatting from E4 (same date style as the rest of the column) before setting the value
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E5").Value = 45063

$excel.CutCopyMode = 0

$ws.Range("E6").Select()
